$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $scratch = $ws.Range("ZZ1")
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
    $scratch.Clear()
}

Set-TextValue 'D2' '66.165.44'
Set-TextValue 'E2' '  +0.29%  '
Set-TextValue 'D3' '3.166.83'
Set-TextValue 'E3' '  -1.48%  '
Set-TextValue 'E4' '  +0.05%  '
Set-TextValue 'D5' '605.19'
Set-TextValue 'E5' '  -0.22%  '
Set-TextValue 'D6' '153.90'
Set-TextValue 'E6' '  +0.24%  '
Set-TextValue 'E7' '  +0.06%  '
Set-TextValue 'D8' '3.165.90'
Set-TextValue 'E8' '  -1.51%  '
Set-TextValue 'D9' '0.547'
Set-TextValue 'E9' '  +2.53%  '
Set-TextValue 'E10' '  -1.49%  '
Set-TextValue 'E11' '  -9.26%  '
Set-TextValue 'D12' '0.517'
Set-TextValue 'E12' '  +1.58%  '
Set-TextValue 'E13' '  -1.63%  '
Set-TextValue 'D14' '38.32'
Set-TextValue 'E14' '  -1.91%  '
Set-TextValue 'D15' '3.688.20'
Set-TextValue 'E15' '  -1.40%  '
Set-TextValue 'D16' '66.198.48'
Set-TextValue 'E16' '  +0.11%  '
Set-TextValue 'D17' '7.41'
Set-TextValue 'E17' '  -0.81%  '
Set-TextValue 'D18' '3.172.92'
Set-TextValue 'E18' '  -1.23%  '
Set-TextValue 'E19' '  +1.02%  '
Set-TextValue 'D20' '509.86'
Set-TextValue 'E20' '  -0.21%  '
Set-TextValue 'E21' '  -0.88%  '
Set-TextValue 'D22' '0.727'
Set-TextValue 'E22' '  -1.16%  '
Set-TextValue 'E23' '  -0.73%  '
Set-TextValue 'D24' '14.72'
Set-TextValue 'E24' '  -3.78%  '
Set-TextValue 'D25' '84.63'
Set-TextValue 'E25' '  -0.81%  '
Set-TextValue 'E26' '  +0.14%  '
Set-TextValue 'E27' '  -0.69%  '
Set-TextValue 'D28' '9.12'
Set-TextValue 'E28' '  -0.28%  '
Set-TextValue 'D29' '2.38'
Set-TextValue 'E29' '  +5.89%  '
Set-TextValue 'D30' '3.04'
Set-TextValue 'E30' '  +6.67%  '
Set-TextValue 'D31' '7.14'
Set-TextValue 'E31' '  +4.67%  '
Set-TextValue 'D32' '27.93'
Set-TextValue 'E32' '  -0.63%  '
Set-TextValue 'E33' '  +0.13%  '
Set-TextValue 'E34' '  -2.43%  '
Set-TextValue 'D35' '6.50'
Set-TextValue 'E35' '  -1.56%  '
Set-TextValue 'D36' '500.47'
Set-TextValue 'E36' '  +4.23%  '
Set-TextValue 'D37' '54.87'
Set-TextValue 'E37' '  -1.00%  '
Set-TextValue 'D38' '0.0881'
Set-TextValue 'E38' '  -2.65%  '
Set-TextValue 'D39' '0.0419'
Set-TextValue 'E39' '  -0.36%  '
Set-TextValue 'E40' '  +7.96%  '
Set-TextValue 'D41' '8.75'
Set-TextValue 'E41' '  -2.77%  '
Set-TextValue 'D42' '0.0₃0681'
Set-TextValue 'E42' '  +6.02%  '
Set-TextValue 'E43' '  -0.34%  '
Set-TextValue 'E44' '  -5.14%  '
Set-TextValue 'E45' '  -1.21%  '
Set-TextValue 'D46' '2.825.04'
Set-TextValue 'E46' '  -4.19%  '
Set-TextValue 'D47' '27.94'
Set-TextValue 'E47' '  -2.53%  '
Set-TextValue 'D48' '2.37'
Set-TextValue 'E48' '  +2.70%  '
Set-TextValue 'E50' '  +0.58%  '
Set-TextValue 'B51' 'Arweave'
Set-TextValue 'C51' 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue 'D51' '35.13'
Set-TextValue 'E51' '  +4.81%  '
